$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 07:52"

# Row 36 - Pakistan: refreshed case numbers
$ws.Range("B36").Value = 4788
$ws.Range("C36").Value = 93
$ws.Range("D36").Value = 762
$ws.Range("E36").Value = 3955
$ws.Range("F36").Value = 50
$ws.Range("G36").Value = 5
$ws.Range("H36").Value = 71

# Rows 48-50 - Tailandia overtakes Catar (and Catar overtakes Colombia) in the ranking
$ws.Range("A48").Value = "Tailandia"
$ws.Range("B48").Value = 2518
$ws.Range("C48").Value = 45
$ws.Range("D48").Value = 1135
$ws.Range("E48").Value = 1348
$ws.Range("F48").Value = 61
$ws.Range("G48").Value = 2
$ws.Range("H48").Value = 35

$ws.Range("A49").Value = "Catar"
$ws.Range("B49").Value = 2512
$ws.Range("C49").Value = 0
$ws.Range("D49").Value = 227
$ws.Range("E49").Value = 2279
$ws.Range("F49").Value = 37
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 6

$ws.Range("A50").Value = "Colombia"
$ws.Range("B50").Value = 2473
$ws.Range("C50").Value = 0
$ws.Range("D50").Value = 197
$ws.Range("E50").Value = 2196
$ws.Range("F50").Value = 85
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 80

# Rows 75-76 - Kazajistan overtakes Camerun in the ranking
$ws.Range("A75").Value = "Kazajistan"
$ws.Range("B75").Value = 840
$ws.Range("C75").Value = 28
$ws.Range("D75").Value = 64
$ws.Range("E75").Value = 766
$ws.Range("F75").Value = 21
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 10

$ws.Range("A76").Value = "Camerun"
$ws.Range("B76").Value = 820
$ws.Range("C76").Value = 0
$ws.Range("D76").Value = 98
$ws.Range("E76").Value = 710
$ws.Range("F76").Value = 0
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 12

# Row 82 - Bulgaria: refreshed case numbers
$ws.Range("B82").Value = 648
$ws.Range("C82").Value = 13
$ws.Range("D82").Value = 62
$ws.Range("E82").Value = 560
$ws.Range("F82").Value = 32
$ws.Range("G82").Value = 1
$ws.Range("H82").Value = 26

# Row 130 - El Salvador: refreshed case numbers
$ws.Range("B130").Value = 118
$ws.Range("C130").Value = 1
$ws.Range("D130").Value = 19
$ws.Range("E130").Value = 93
